# Apply odds updates to Sheet1 (Jogos_da_Semana_FlashScore_2024-10-13.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("G6").Value = 2.2
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 10.75
$ws.Range("Z6").Value = 22
$ws.Range("AC6").Value = 9.25
$ws.Range("AD6").Value = 6.2
$ws.Range("AF6").Value = 65
$ws.Range("AK6").Value = 27
$ws.Range("AP6").Value = 18.5
$ws.Range("AU6").Value = 6.7
$ws.Range("AW6").Value = 5

# Row 12
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.4
$ws.Range("Q12").Value = 2.03
$ws.Range("R12").Value = 1.83

# Row 13
$ws.Range("G13").Value = 2.05
$ws.Range("I13").Value = 3.5
$ws.Range("L13").Value = 4
$ws.Range("X13").Value = 9.5
$ws.Range("Z13").Value = 19
$ws.Range("AA13").Value = 19
$ws.Range("AH13").Value = 17
$ws.Range("AK13").Value = 29
$ws.Range("AO13").Value = 12
$ws.Range("AR13").Value = 67

# Row 19
$ws.Range("Q19").Value = 1.98
$ws.Range("R19").Value = 1.88
